# cap nhat thong ke
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 13 (Nhom N12): set "Diem bt hackerrank" to "400/400" and "Diem tong ket" to 8
$ws.Range("B13").Value = "400/400"
$ws.Range("D13").Value = 8

# Row 19 (Nhom N18): set "Diem bt hackerrank" to "400/400" and "Diem tong ket" to 8
$ws.Range("B19").Value = "400/400"
$ws.Range("D19").Value = 8

# Update the view: scroll back to top-left A1 and select D19
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D19").Select()
